$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 5.1
$ws.Range("D1").Value = 4.3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 4
$ws.Range("G1").Value = 4
$ws.Range("H1").Value = 4

$ws.Range("C2").Value = 7.4
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 5

$ws.Range("C3").Value = 7.5
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 6
$ws.Range("G3").Value = 6

$ws.Range("C4").Value = 11.6
$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 7

$ws.Range("C5").Value = 8.800000000000001
$ws.Range("D5").Value = 8
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 8
$ws.Range("G5").Value = 8
$ws.Range("H5").Value = 8

$ws.Range("C6").Value = 12.4
$ws.Range("D6").Value = 9
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 9
$ws.Range("G6").Value = 9
$ws.Range("H6").Value = 9

$ws.Range("C7").Value = 13.1
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 10
$ws.Range("G7").Value = 10
$ws.Range("H7").Value = 10

$ws.Range("C8").Value = 12.4
$ws.Range("D8").Value = 11.7
$ws.Range("E8").Value = 11
$ws.Range("F8").Value = 11
$ws.Range("G8").Value = 11
$ws.Range("H8").Value = 11

$ws.Range("C9").Value = 13.6
$ws.Range("D9").Value = 12
$ws.Range("E9").Value = 12
$ws.Range("F9").Value = 12
$ws.Range("G9").Value = 12
$ws.Range("H9").Value = 12
